# Issue #5: fund sheet (基金受益憑證 / sheet4) — fill in header row + append the
# trailing legislator/meta columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index), matching the pattern
# already used on the other sheets (土地/建物/存款 etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# ---- Row 1: turn the accidental duplicate data row into a proper header ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# Give the new header cells (I1:O1) the same bold/border/centered style as
# the existing header cells (B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1:O1").PasteSpecial(-4122)

# ---- Rows 2-12: append the trailing metadata columns ----
$lastRow = 12
for ($r = 2; $r -le $lastRow; $r++) {
    $idx = $ws.Cells.Item($r, 1).Value
    $ws.Cells.Item($r, 9).Value  = "fund"
    $ws.Cells.Item($r, 10).Value = "normal"
    $ws.Cells.Item($r, 11).Value = "2012-04-09"
    $ws.Cells.Item($r, 12).Value = "吳秉叡"
    $ws.Cells.Item($r, 13).Value = 1324
    $ws.Cells.Item($r, 14).Value = "tmp8f8d1"
    $ws.Cells.Item($r, 15).Value = $idx
}
